# Auto-generated edit script: updates cryptos list values per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.000.15"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.561.42"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.490"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0599"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0854"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "1.783.61"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "1.528.78"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "26.984.41"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  +1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("E31").Value = "  +2.50%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("D34").Value = "1.426.20"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("E36").Value = "  +8.26%  "
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "1.699.11"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  +4.80%  "
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0961"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
